# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.295.87'
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = '3.391.14'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.39%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("E9").Value = '  +7.52%  '
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.60'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("E12").Value = '  +3.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '679.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '3.931.07'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").Value = '69.380.61'
$ws.Range("E16").Value = '  +2.32%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.394.00'
$ws.Range("E17").Value = '  +2.01%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.120'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.30%  '
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.81%  '
$ws.Range("E29").Value = '  +2.25%  '
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '557.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.28%  '
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '3.677.14'
$ws.Range("E37").Value = '  -0.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("E39").Value = '  +3.55%  '
$ws.Range("D40").Value = '0.0₃0719'
$ws.Range("E40").Value = '  +6.55%  '
$ws.Range("E41").Value = '  +3.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("E43").Value = '  +1.41%  '
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.29%  '

"done"
